$wb = $excel.ActiveWorkbook

# Insert a new "Year_selection" worksheet right after "Technology_selection"
# (i.e. before "Technology_in_region_selection").
$afterSheet = $wb.Worksheets.Item("Technology_selection")
$yearSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$yearSheet.Name = "Year_selection"

# Header row
$yearSheet.Range("A1").Value = "Year"
$yearSheet.Range("B1").Value = "Year selection"

# Data rows
$years = @(2015, 2020, 2025, 2030, 2035, 2040, 2045, 2050)
$selected = @(1, 1, 0, 1, 1, 0, 1, 1)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $yearSheet.Cells.Item($row, 1).Value = $years[$i]
    $yearSheet.Cells.Item($row, 2).Value = $selected[$i]
}

# Match the selection/active-cell state recorded in the target file.
$yearSheet.Range("B10").Select()

# Make the new sheet the active one (matches activeTab pointing at it).
$yearSheet.Activate()
